$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 3676821.2
$ws.Range("J9").Value = 302
$ws.Range("L9").Value = 302
$ws.Range("N9").Value = -640
$ws.Range("H17").Value = 2269.111
$ws.Range("J17").Value = 2269.111
$ws.Range("L17").Value = 6807.333
$ws.Range("N17").Value = -7143.333
$ws.Range("H33").Value = 164.81818
$ws.Range("I33").Value = 164.81818
$ws.Range("K33").Value = 164.81818
$ws.Range("M33").Value = 64.18181999999999
$ws.Range("H62").Value = 5811
$ws.Range("I62").Value = 5811
$ws.Range("K62").Value = 5811
$ws.Range("M62").Value = -5187
$ws.Range("H65").Value = 5811
$ws.Range("I65").Value = 5811
$ws.Range("K65").Value = 29055
$ws.Range("M65").Value = -25935
$ws.Range("H111").Value = 1978.6
$ws.Range("J111").Value = 1396.5
$ws.Range("L111").Value = 4189.5
$ws.Range("N111").Value = -10323.5
$ws.Range("H138").Value = 324471.38
$ws.Range("I138").Value = 4196.36
$ws.Range("K138").Value = 12589.08
$ws.Range("M138").Value = -7449.079999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6093.431
$ws.Range("I32").Value = 4971.3257
$ws.Range("K32").Value = 4971.3257
$ws.Range("M32").Value = -4684.3257
$ws.Range("H61").Value = 3996.074
$ws.Range("I61").Value = 1942.5
$ws.Range("K61").Value = 1942.5
$ws.Range("M61").Value = -1730.5
$ws.Range("H97").Value = 944
$ws.Range("I97").Value = 1024.3
$ws.Range("K97").Value = 1024.3
$ws.Range("M97").Value = -528.3
$ws.Range("H136").Value = 3996.074
$ws.Range("I136").Value = 1942.5
$ws.Range("K136").Value = 5827.5
$ws.Range("M136").Value = -3277.5
$ws.Range("H139").Value = 74152.30499999999
$ws.Range("J139").Value = 74152.30499999999
$ws.Range("L139").Value = 74152.30499999999
$ws.Range("N139").Value = -84432.30499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 23813768
$ws.Range("J20").Value = 2699
$ws.Range("L20").Value = 2699
$ws.Range("N20").Value = -3193
$ws.Range("H94").Value = 52632148
$ws.Range("I94").Value = 52632148
$ws.Range("K94").Value = 52632148
$ws.Range("M94").Value = -52631697
$ws.Range("H134").Value = 3153.6206
$ws.Range("I134").Value = 2401.9524
$ws.Range("K134").Value = 7205.8572
$ws.Range("M134").Value = -4670.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1849.6666
$ws.Range("I16").Value = 1849.6666
$ws.Range("K16").Value = 1849.6666
$ws.Range("M16").Value = -1562.6666
$ws.Range("H31").Value = 3571.1887
$ws.Range("I31").Value = 2982.7
$ws.Range("K31").Value = 2982.7
$ws.Range("M31").Value = -2687.7
$ws.Range("H34").Value = 3571.1887
$ws.Range("I34").Value = 2982.7
$ws.Range("K34").Value = 2982.7
$ws.Range("M34").Value = -2780.7
$ws.Range("H113").Value = 1849.6666
$ws.Range("I113").Value = 1849.6666
$ws.Range("K113").Value = 1849.6666
$ws.Range("M113").Value = 320.3334
$ws.Range("H132").Value = 13892327
$ws.Range("J132").Value = 4999.857
$ws.Range("L132").Value = 14999.571
$ws.Range("N132").Value = -20059.571
$ws.Range("H134").Value = 2867.5833
$ws.Range("I134").Value = 2312.3157
$ws.Range("J134").Value = 4977.6
$ws.Range("K134").Value = 6936.9471
$ws.Range("L134").Value = 14932.8
$ws.Range("M134").Value = -4401.9471
$ws.Range("N134").Value = -20002.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 237.11111
$ws.Range("J86").Value = 184
$ws.Range("L86").Value = 552
$ws.Range("N86").Value = -2924
$ws.Range("H89").Value = 237.11111
$ws.Range("J89").Value = 184
$ws.Range("L89").Value = 1656
$ws.Range("N89").Value = -13512
$ws.Range("H115").Value = 431357
$ws.Range("J115").Value = 602999.8
$ws.Range("L115").Value = 1808999.4
$ws.Range("N115").Value = -1811349.4
$ws.Range("H131").Value = 8630.375
$ws.Range("I131").Value = 16701.857
$ws.Range("J131").Value = 2352.5557
$ws.Range("K131").Value = 50105.571
$ws.Range("L131").Value = 7057.6671
$ws.Range("M131").Value = -45065.571
$ws.Range("N131").Value = -17137.6671
$ws.Range("H132").Value = 1971.8334
$ws.Range("J132").Value = 2537.4
$ws.Range("L132").Value = 22836.6
$ws.Range("N132").Value = -27896.6
$ws.Range("H139").Value = 6689.391
$ws.Range("I139").Value = 7928.625
$ws.Range("J139").Value = 3856.8572
$ws.Range("K139").Value = 23785.875
$ws.Range("L139").Value = 11570.5716
$ws.Range("M139").Value = -18645.875
$ws.Range("N139").Value = -21850.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 37038624
$ws.Range("J80").Value = 1828.7142
$ws.Range("L80").Value = 1828.7142
$ws.Range("N80").Value = -3824.7142
$ws.Range("H83").Value = 37038624
$ws.Range("J83").Value = 1828.7142
$ws.Range("L83").Value = 9143.571
$ws.Range("N83").Value = -19127.571
$ws.Range("H113").Value = 4712.0386
$ws.Range("I113").Value = 4795.381
$ws.Range("J113").Value = 4362
$ws.Range("K113").Value = 4795.381
$ws.Range("L113").Value = 4362
$ws.Range("M113").Value = -2625.381
$ws.Range("N113").Value = -8702
$ws.Range("H132").Value = 2625.743
$ws.Range("I132").Value = 2785.1924
$ws.Range("J132").Value = 2165.111
$ws.Range("K132").Value = 8355.5772
$ws.Range("L132").Value = 6495.333
$ws.Range("M132").Value = -5825.5772
$ws.Range("N132").Value = -11555.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4635.196
$ws.Range("I40").Value = 4902.6387
$ws.Range("J40").Value = 3672.4
$ws.Range("K40").Value = 4902.6387
$ws.Range("L40").Value = 3672.4
$ws.Range("M40").Value = -4766.6387
$ws.Range("N40").Value = -3944.4
$ws.Range("H61").Value = 1272.25
$ws.Range("I61").Value = 1149.5834
$ws.Range("J61").Value = 1640.25
$ws.Range("K61").Value = 1149.5834
$ws.Range("L61").Value = 1640.25
$ws.Range("M61").Value = -947.5834
$ws.Range("N61").Value = -2044.25
$ws.Range("H100").Value = 1129366.2
$ws.Range("I100").Value = 1354574.6
$ws.Range("J100").Value = 3324.75
$ws.Range("K100").Value = 1354574.6
$ws.Range("L100").Value = 3324.75
$ws.Range("M100").Value = -1354033.6
$ws.Range("N100").Value = -4406.75
$ws.Range("H113").Value = 1272.25
$ws.Range("I113").Value = 1149.5834
$ws.Range("J113").Value = 1640.25
$ws.Range("K113").Value = 1149.5834
$ws.Range("L113").Value = 1640.25
$ws.Range("M113").Value = 1020.4166
$ws.Range("N113").Value = -5980.25
$ws.Range("H122").Value = 3364.6316
$ws.Range("I122").Value = 3545.1428
$ws.Range("K122").Value = 10635.4284
$ws.Range("M122").Value = -8185.428400000001
$ws.Range("H132").Value = 3031.7126
$ws.Range("I132").Value = 2357.234
$ws.Range("K132").Value = 7071.701999999999
$ws.Range("M132").Value = -4541.701999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 697.8125
$ws.Range("I107").Value = 801.25
$ws.Range("J107").Value = 387.5
$ws.Range("K107").Value = 2403.75
$ws.Range("L107").Value = 1162.5
$ws.Range("M107").Value = -483.75
$ws.Range("N107").Value = -5002.5
$ws.Range("H132").Value = 6538375.5
$ws.Range("I132").Value = 7939048
$ws.Range("K132").Value = 23817144
$ws.Range("M132").Value = -23814614
